$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 40 & 42: Coin/Link swapped (ARBITRUM <-> Monero), with refreshed Price/Volume(1h) ---
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$d40 = $ws.Range("D40")
$d40.NumberFormat = "@"
$d40.Value = "137.55"
$ws.Range("E40").Value = "  -0.19%  "

$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$d42 = $ws.Range("D42")
$d42.NumberFormat = "@"
$d42.Value = "1.98"
$ws.Range("E42").Value = "  -0.81%  "

# --- Refresh Price (D) / Volume(1h) (E) figures for every other listed row ---
# NumberFormat is forced to Text ("@") before writing Price values so that
# numeric-looking strings (e.g. "408.00", "1.00") are preserved exactly,
# matching how the source data keeps these as literal text.
$d = $ws.Range("D2")
$d.NumberFormat = "@"
$d.Value = "62.086.89"
$ws.Range("E2").Value = "  -1.31%  "
$d = $ws.Range("D3")
$d.NumberFormat = "@"
$d.Value = "3.410.76"
$ws.Range("E3").Value = "  -1.88%  "
$ws.Range("E4").Value = "  -0.02%  "
$d = $ws.Range("D5")
$d.NumberFormat = "@"
$d.Value = "408.00"
$ws.Range("E5").Value = "  -0.89%  "
$d = $ws.Range("D6")
$d.NumberFormat = "@"
$d.Value = "134.33"
$ws.Range("E6").Value = "  +4.14%  "
$d = $ws.Range("D7")
$d.NumberFormat = "@"
$d.Value = "0.592"
$ws.Range("E7").Value = "  -1.25%  "
$d = $ws.Range("D8")
$d.NumberFormat = "@"
$d.Value = "1.00"
$ws.Range("E8").Value = "  -0.03%  "
$d = $ws.Range("D9")
$d.NumberFormat = "@"
$d.Value = "0.688"
$ws.Range("E9").Value = "  -2.18%  "
$d = $ws.Range("D10")
$d.NumberFormat = "@"
$d.Value = "0.122"
$ws.Range("E10").Value = "  -5.87%  "
$d = $ws.Range("D11")
$d.NumberFormat = "@"
$d.Value = "42.68"
$ws.Range("E11").Value = "  -1.48%  "
$ws.Range("E12").Value = "  -1.05%  "
$d = $ws.Range("D13")
$d.NumberFormat = "@"
$d.Value = "8.43"
$ws.Range("E13").Value = "  -3.74%  "
$d = $ws.Range("D14")
$d.NumberFormat = "@"
$d.Value = "19.89"
$ws.Range("E14").Value = "  -1.57%  "
$d = $ws.Range("D15")
$d.NumberFormat = "@"
$d.Value = "3.413.21"
$ws.Range("E15").Value = "  -1.45%  "
$d = $ws.Range("D16")
$d.NumberFormat = "@"
$d.Value = "62.066.96"
$ws.Range("E16").Value = "  -1.12%  "
$ws.Range("E17").Value = "  -2.51%  "
$d = $ws.Range("D18")
$d.NumberFormat = "@"
$d.Value = "11.03"
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("E19").Value = "  -4.79%  "
$ws.Range("E20").Value = "  -4.65%  "
$d = $ws.Range("D21")
$d.NumberFormat = "@"
$d.Value = "84.29"
$ws.Range("E21").Value = "  +2.15%  "
$d = $ws.Range("D22")
$d.NumberFormat = "@"
$d.Value = "314.36"
$ws.Range("E22").Value = "  -0.11%  "
$d = $ws.Range("D23")
$d.NumberFormat = "@"
$d.Value = "12.89"
$ws.Range("E23").Value = "  -2.26%  "
$d = $ws.Range("D24")
$d.NumberFormat = "@"
$d.Value = "3.16"
$ws.Range("E24").Value = "  -0.50%  "
$d = $ws.Range("D25")
$d.NumberFormat = "@"
$d.Value = "4.77"
$ws.Range("E25").Value = "  +9.29%  "
$d = $ws.Range("D26")
$d.NumberFormat = "@"
$d.Value = "29.64"
$ws.Range("E26").Value = "  -2.65%  "
$d = $ws.Range("D27")
$d.NumberFormat = "@"
$d.Value = "8.21"
$ws.Range("E27").Value = "  -0.02%  "
$d = $ws.Range("D28")
$d.NumberFormat = "@"
$d.Value = "2.85"
$ws.Range("E28").Value = "  +6.09%  "
$ws.Range("E29").Value = "  -2.44%  "
$d = $ws.Range("D30")
$d.NumberFormat = "@"
$d.Value = "0.175"
$ws.Range("E30").Value = "  -3.48%  "
$ws.Range("E31").Value = "  -2.93%  "
$d = $ws.Range("D32")
$d.NumberFormat = "@"
$d.Value = "42.71"
$ws.Range("E32").Value = "  -3.46%  "
$ws.Range("E33").Value = "  -0.11%  "
$d = $ws.Range("D34")
$d.NumberFormat = "@"
$d.Value = "11.37"
$ws.Range("E34").Value = "  -6.32%  "
$d = $ws.Range("D35")
$d.NumberFormat = "@"
$d.Value = "0.0483"
$ws.Range("E35").Value = "  -2.06%  "
$d = $ws.Range("D36")
$d.NumberFormat = "@"
$d.Value = "51.74"
$ws.Range("E36").Value = "  -1.84%  "
$ws.Range("E37").Value = "  +0.23%  "
$d = $ws.Range("D38")
$d.NumberFormat = "@"
$d.Value = "3.41"
$ws.Range("E38").Value = "  -4.76%  "
$d = $ws.Range("D39")
$d.NumberFormat = "@"
$d.Value = "2.95"
$ws.Range("E39").Value = "  -2.91%  "
$ws.Range("E41").Value = "  -0.41%  "
$ws.Range("E43").Value = "  +4.06%  "
$d = $ws.Range("D44")
$d.NumberFormat = "@"
$d.Value = "4.02"
$ws.Range("E44").Value = "  +1.13%  "
$d = $ws.Range("D45")
$d.NumberFormat = "@"
$d.Value = "16.75"
$ws.Range("E45").Value = "  -5.88%  "
$ws.Range("E46").Value = "  -2.19%  "
$d = $ws.Range("D47")
$d.NumberFormat = "@"
$d.Value = "21.35"
$ws.Range("E47").Value = "  -4.78%  "
$d = $ws.Range("D48")
$d.NumberFormat = "@"
$d.Value = "2.122.88"
$ws.Range("E48").Value = "  -4.42%  "
$d = $ws.Range("D49")
$d.NumberFormat = "@"
$d.Value = "2.31"
$ws.Range("E49").Value = "  -2.86%  "
$d = $ws.Range("D50")
$d.NumberFormat = "@"
$d.Value = "1.91"
$ws.Range("E50").Value = "  +0.52%  "
$d = $ws.Range("D51")
$d.NumberFormat = "@"
$d.Value = "1.68"
$ws.Range("E51").Value = "  +19.81%  "
